$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '49.781.76'
$ws.Range('E2').Value = '  +2.86%  '

$ws.Range('D3').Value = '2.620.46'
$ws.Range('E3').Value = '  +4.54%  '

$ws.Range('E4').Value = '  -0.09%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '328.46'
$ws.Range('E5').Value = '  +2.11%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '111.00'
$ws.Range('E6').Value = '  +2.32%  '

$ws.Range('E7').Value = '  +1.10%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  -0.01%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.559'
$ws.Range('E9').Value = '  +3.24%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.90'

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.73'
$ws.Range('E11').Value = '  +2.00%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0821'
$ws.Range('E12').Value = '  +0.15%  '

$ws.Range('E13').Value = '  +0.66%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.29'
$ws.Range('E14').Value = '  +1.27%  '

$ws.Range('D15').Value = '3.035.98'
$ws.Range('E15').Value = '  +4.73%  '

$ws.Range('D16').Value = '2.622.26'
$ws.Range('E16').Value = '  +4.67%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.875'
$ws.Range('E17').Value = '  +3.55%  '

$ws.Range('D18').Value = '49.728.86'
$ws.Range('E18').Value = '  +3.11%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.10'
$ws.Range('E19').Value = '  +10.61%  '

$ws.Range('E20').Value = '  +1.82%  '

$ws.Range('E21').Value = '  +0.42%  '

$ws.Range('D22').Value = '0.0₃0954'
$ws.Range('E22').Value = '  +0.76%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '281.40'
$ws.Range('E23').Value = '  +0.34%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '72.79'
$ws.Range('E24').Value = '  +0.56%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.59'
$ws.Range('E25').Value = '  +1.22%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.62'
$ws.Range('E26').Value = '  +3.20%  '

$ws.Range('E27').Value = '  -0.05%  '

$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.01'
$ws.Range('E28').Value = '  +1.97%  '

$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.23'
$ws.Range('E29').Value = '  -2.54%  '

$ws.Range('E30').Value = '  +1.90%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '36.22'
$ws.Range('E31').Value = '  +2.29%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '49.80'
$ws.Range('E32').Value = '  +0.84%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.80'
$ws.Range('E33').Value = '  +0.58%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.45'
$ws.Range('E34').Value = '  +1.72%  '

$ws.Range('E35').Value = '  -0.26%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0794'
$ws.Range('E36').Value = '  +1.17%  '

$ws.Range('E37').Value = '  +5.23%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.74'
$ws.Range('E38').Value = '  +1.76%  '

$ws.Range('E39').Value = '  +5.68%  '

$ws.Range('B40').Value = 'Monero'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '124.09'
$ws.Range('E40').Value = '  +1.97%  '

$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '22.73'
$ws.Range('E41').Value = '  +5.59%  '

$ws.Range('B42').Value = 'Stellar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.112'
$ws.Range('E42').Value = '  +0.63%  '

$ws.Range('E43').Value = '  +0.36%  '

$ws.Range('E44').Value = '  +4.13%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.36'
$ws.Range('E45').Value = '  +5.72%  '

$ws.Range('D46').Value = '2.054.70'
$ws.Range('E46').Value = '  +2.10%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.23'
$ws.Range('E47').Value = '  +11.92%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.02'
$ws.Range('E48').Value = '  +8.93%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.06'
$ws.Range('E49').Value = '  +0.27%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.37'
$ws.Range('E50').Value = '  +3.54%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '81.98'
$ws.Range('E51').Value = '  +1.64%  '
